# Update the "Förändrad" (Changed) date column (column C) from
# 2023-09-19 (Excel serial 45188) to 2023-09-20 (Excel serial 45189)
# for every data row (rows 2-203) in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C203").Value = 45189
